$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header cells (bold, matching the rest of row 1) ------
$ws.Range("H1").Value = "Guild"
$ws.Range("I1").Value = "Species_of_concern"
$ws.Range("H1:I1").Font.Bold = $true

# --- Per-species Guild ("H") and Species_of_concern ("I") values ----------
# Keyed by row number (row 1 is the header; data runs 2-28).
$guild = @{
   2 = "Bird"         # Russet-naped wood rail
   3 = "Frugivore"     # Red brocket
   4 = "Ungulate"      # White-tailed deer
   5 = "Ungulate"      # Collared peccary
   6 = "Ungulate"      # White-lipped peccary
   7 = "Carnivore"     # Gray fox
   8 = "Carnivore"     # Jaguar
   9 = "Carnivore"     # Ocelot
  10 = "Carnivore"     # Puma
  11 = "Carnivore"     # White-nosed coati
  12 = "Carnivore"     # Margay
  13 = "Carnivore"     # Tayra
  14 = "Carnivore"     # Striped hog-nosed skunk
  15 = "Insectivore"   # Nine-banded armadillo
  16 = "Bird"          # Gray-headed dove
  17 = "Bird"          # Ruddy quail-dove
  18 = "Generalist"    # Common opossum
  19 = "Insectivore"   # Brown four-eyed opossum
  20 = "Bird"          # Crested guan
  21 = "Bird"          # Great curassow
  22 = "Bird"          # Ocellated turkey
  23 = "Bird"          # Wood thrush
  24 = "Bird"          # Bare-throated tiger heron
  25 = "Ungulate"      # Baird's tapir
  26 = "Generalist"    # Agouti
  27 = "Frugivore"     # Paca
  28 = "Bird"          # Great tinamou
}

$concern = @{
   2 = "N"
   3 = "Y"
   4 = "N"
   5 = "N"
   6 = "Y"
   7 = "N"
   8 = "Y"
   9 = "Y"
  10 = "Y"
  11 = "N"
  12 = "Y"
  13 = "N"
  14 = "N"
  15 = "N"
  16 = "N"
  17 = "N"
  18 = "N"
  19 = "N"
  20 = "Y"
  21 = "Y"
  22 = "Y"
  23 = "N"
  24 = "N"
  25 = "Y"
  26 = "N"
  27 = "N"
  28 = "N"
}

foreach ($r in ($guild.Keys | Sort-Object)) {
    $ws.Cells.Item($r, 8).Value = $guild[$r]
    $ws.Cells.Item($r, 9).Value = $concern[$r]
}

# --- Drop the leftover highlight style from rows 18 & 27 -------------------
# (the source workbook no longer marks these rows/cells with the special
# "applyFill" style once the new columns are in place)
$ws.Rows.Item(18).ClearFormats()
$ws.Range("A27:B27").ClearFormats()

# --- Selection / printing touch-ups matching the refreshed workbook -------
$ws.Range("I29").Select() | Out-Null
$ws.PageSetup.Orientation = 1 | Out-Null
